# Insert a new data row at row 24 (pushing the existing rows 24-102 down to
# 25-103) and populate the newly inserted row with the new reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("24:24").Insert()

$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value2 = 44883
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100101
$ws.Range("H24").Value = "Berries"
$ws.Range("I24").Value = 100101001
$ws.Range("J24").Value = "Arándano (blue)"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 3200
$ws.Range("O24").Value = 3200
$ws.Range("P24").Value = 3200
$ws.Range("Q24").Value = '$/kilo'
$ws.Range("R24").Value = "Región del Maule"
$ws.Range("S24").Value = 3200
$ws.Range("T24").Value = 1
